$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Resize column A (user manually narrowed it)
$ws.Columns.Item(1).ColumnWidth = 39.1868

# 2. Style + fill the previously-empty F column cells (rows 2-9, 11, 12) with
#    centered alignment (same look as the rest of the table)
$fRows = @(2,3,4,5,6,7,8,9,11,12)
foreach ($r in $fRows) {
    $cell = $ws.Cells.Item($r, 6)
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4108
}

# 3. Fill in the rest of row 10 (task finished ahead of schedule)
$d10 = $ws.Range("D10")
$d10.HorizontalAlignment = -4108
$d10.VerticalAlignment = -4108
$d10.Value = "entregue dentro do prazo"

$e10 = $ws.Range("E10")
$e10.HorizontalAlignment = -4108
$e10.VerticalAlignment = -4108
$e10.Value = "finalizado"

$f10 = $ws.Range("F10")
$f10.HorizontalAlignment = -4108
$f10.VerticalAlignment = -4108
$f10.Value = 43204
$f10.NumberFormat = "DD/MM/YY"

# 4. Row heights grow to fit the newly wrapped text after the column resize
$ws.Rows.Item(3).RowHeight = 53.95
$ws.Rows.Item(4).RowHeight = 50.95
$ws.Rows.Item(10).RowHeight = 47.2

# 5. Leave the selection on A10
$ws.Range("A10").Select()
